$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.65  # G2
$ws.Cells.Item(2, 8).Value = 5.6  # H2
$ws.Cells.Item(2, 9).Value = 6.6  # I2
$ws.Cells.Item(2, 10).Value = 4.4  # J2
$ws.Cells.Item(2, 15).Value = 1.28  # O2
$ws.Cells.Item(2, 16).Value = 2.06  # P2
$ws.Cells.Item(2, 17).Value = 1.76  # Q2
$ws.Cells.Item(2, 18).Value = 1.42  # R2
$ws.Cells.Item(2, 20).Value = 1.9  # T2
$ws.Cells.Item(2, 23).Value = 2.52  # W2
$ws.Cells.Item(2, 24).Value = 25  # X2
$ws.Cells.Item(2, 25).Value = 40  # Y2
$ws.Cells.Item(2, 28).Value = 16  # AB2
$ws.Cells.Item(2, 29).Value = 14.5  # AC2
$ws.Cells.Item(2, 30).Value = 55  # AD2
$ws.Cells.Item(2, 32).Value = 40  # AF2
$ws.Cells.Item(2, 33).Value = 40  # AG2
$ws.Cells.Item(2, 36).Value = 180  # AJ2
$ws.Cells.Item(2, 37).Value = 65  # AK2
$ws.Cells.Item(2, 40).Value = 10  # AN2
$ws.Cells.Item(3, 6).Value = 7.6  # F3
$ws.Cells.Item(3, 7).Value = 7.8  # G3
$ws.Cells.Item(3, 8).Value = 1.42  # H3
$ws.Cells.Item(3, 10).Value = 5.3  # J3
$ws.Cells.Item(3, 11).Value = 5.9  # K3
$ws.Cells.Item(3, 12).Value = 1.29  # L3
$ws.Cells.Item(3, 14).Value = 5.2  # N3
$ws.Cells.Item(3, 16).Value = 2.46  # P3
$ws.Cells.Item(3, 17).Value = 1.57  # Q3
$ws.Cells.Item(3, 18).Value = 1.57  # R3
$ws.Cells.Item(3, 20).Value = 1.85  # T3
$ws.Cells.Item(3, 21).Value = 2  # U3
$ws.Cells.Item(3, 27).Value = 12.5  # AA3
$ws.Cells.Item(3, 30).Value = 10.5  # AD3
$ws.Cells.Item(3, 33).Value = 36  # AG3
$ws.Cells.Item(3, 34).Value = 990  # AH3
$ws.Cells.Item(3, 35).Value = 1000  # AI3
$ws.Cells.Item(4, 7).Value = 1.53  # G4
$ws.Cells.Item(4, 8).Value = 3.85  # H4
$ws.Cells.Item(4, 9).Value = 42  # I4
$ws.Cells.Item(4, 10).Value = 3.15  # J4
$ws.Cells.Item(4, 12).Value = 1.02  # L4
$ws.Cells.Item(4, 13).Value = 1.05  # M4
$ws.Cells.Item(4, 14).Value = 1.1  # N4
$ws.Cells.Item(4, 16).Value = 1.53  # P4
$ws.Cells.Item(4, 17).Value = 1.2  # Q4
$ws.Cells.Item(4, 18).Value = 1.22  # R4
$ws.Cells.Item(4, 19).Value = 1.66  # S4
$ws.Cells.Item(4, 23).Value = 2.72  # W4
$ws.Cells.Item(4, 40).Value = 29  # AN4
$ws.Cells.Item(5, 6).Value = 2.18  # F5
$ws.Cells.Item(5, 10).Value = 3.2  # J5
$ws.Cells.Item(5, 17).Value = 2.12  # Q5
$ws.Cells.Item(5, 18).Value = 1.27  # R5
$ws.Cells.Item(5, 22).Value = 1.33  # V5
$ws.Cells.Item(5, 23).Value = 1.72  # W5
$ws.Cells.Item(5, 27).Value = 900  # AA5
$ws.Cells.Item(5, 31).Value = 150  # AE5
$ws.Cells.Item(5, 33).Value = 12  # AG5
$ws.Cells.Item(5, 41).Value = 1000  # AO5
$ws.Cells.Item(6, 6).Value = 1.77  # F6
$ws.Cells.Item(6, 7).Value = 1.9  # G6
$ws.Cells.Item(6, 11).Value = 3.75  # K6
$ws.Cells.Item(6, 14).Value = 2.92  # N6
$ws.Cells.Item(6, 15).Value = 1.42  # O6
$ws.Cells.Item(6, 16).Value = 1.64  # P6
$ws.Cells.Item(6, 17).Value = 2.26  # Q6
$ws.Cells.Item(6, 20).Value = 2.04  # T6
$ws.Cells.Item(6, 23).Value = 2.1  # W6
$ws.Cells.Item(6, 26).Value = 130  # Z6
$ws.Cells.Item(6, 30).Value = 27  # AD6
$ws.Cells.Item(6, 34).Value = 60  # AH6
$ws.Cells.Item(6, 36).Value = 60  # AJ6
$ws.Cells.Item(6, 37).Value = 27  # AK6
$ws.Cells.Item(6, 38).Value = 290  # AL6
